# Update the "addCourse_details" sheet: change the category test data
# used when running the Add-Course test from "Teamcity" to "Java", and
# move the active selection to J7 (matches updated utility / grid run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("addCourse_details")

# Update the CategoryName value cell (H2) from "Teamcity" to "Java"
$ws.Range("H2").Value = "Java"

# Reflect the new cursor/selection position left after editing
$ws.Range("J7").Select()
